# Weekly driver report update for 2025-04-19
# Rewrites the "Bad Drivers" and "Good Drivers" tables on the single
# "Driver Summary" sheet with refreshed sample data, growing both tables
# by a few rows and widening column A slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Copy cell FORMATTING (not values) from rows that already carry the
#    styles we need into the rows that are brand new / need a style they
#    don't currently have. Do this before any value is overwritten so the
#    "donor" cells are still pristine.
# ---------------------------------------------------------------------

# New "Bad Drivers" data row (row 5) needs the same look as the existing
# data rows (row 3 / row 4): no border, right aligned numbers.
$ws.Range("A3:D3").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Row 6 becomes the new "Totals:" row -> copy the old Totals formatting
# (currently on row 5: bold label, bold #,##0 numbers) down into row 6.
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)

# The "Good Drivers" section header block shifts from rows 11/12 down to
# rows 12/13. Copy row 12's (column header) formatting down to row 13
# FIRST, while row 12 still has its original look, then stamp row 11's
# (section title) formatting onto row 12.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# Four brand new "Good Drivers" data rows (19-22) need the same styling
# as the existing data rows, e.g. row 14 (plain label, #,##0 sample
# count, right-aligned percentage / vintage).
$ws.Range("A14:E14").Copy()
$ws.Range("A19:E22").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Widen column A (44 -> 45 characters).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 44.16

# ---------------------------------------------------------------------
# 3) Touch J27 (format only, value stays empty) so the sheet's used
#    range / dimension grows to A1:J27, matching the padded blank rows
#    and columns already reserved via column widths F:J.
# ---------------------------------------------------------------------
$ws.Cells.Item(27, 10).Borders.Item(7).LineStyle = -4142

# ---------------------------------------------------------------------
# 4) Write the refreshed "Bad Drivers" table values.
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value2 = "iwlwifi"
$ws.Cells.Item(3, 2).Value2 = 7
$ws.Cells.Item(3, 3).Value2 = 4238
$ws.Cells.Item(3, 4).Value2 = 91.59999999999999

$ws.Cells.Item(4, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4"
$ws.Cells.Item(4, 2).Value2 = 9
$ws.Cells.Item(4, 3).Value2 = 644
$ws.Cells.Item(4, 4).Value2 = 96.7

$ws.Cells.Item(5, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1"
$ws.Cells.Item(5, 2).Value2 = 7
$ws.Cells.Item(5, 3).Value2 = 434
$ws.Cells.Item(5, 4).Value2 = 98.8

$ws.Cells.Item(6, 1).Value2 = "Totals:"
$ws.Cells.Item(6, 2).Value2 = 23
$ws.Cells.Item(6, 3).Value2 = 5316

# ---------------------------------------------------------------------
# 5) Write the refreshed "Good Drivers (Roaming > 99.8%)" table values.
# ---------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value2 = "Good Drivers (Roaming > 99.8%)"

$ws.Cells.Item(13, 1).Value2 = "Adapter-Driver"
$ws.Cells.Item(13, 2).Value2 = "Total Samples"
$ws.Cells.Item(13, 3).Value2 = ""
$ws.Cells.Item(13, 4).Value2 = "Good Roaming Calculation (%)"
$ws.Cells.Item(13, 5).Value2 = "Driver Vintage"

$ws.Cells.Item(14, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Cells.Item(14, 2).Value2 = 10661
$ws.Cells.Item(14, 3).Value2 = ""
$ws.Cells.Item(14, 4).Value2 = 100
$ws.Cells.Item(14, 5).ClearContents()

$ws.Cells.Item(15, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Cells.Item(15, 2).Value2 = 56018
$ws.Cells.Item(15, 3).Value2 = ""
$ws.Cells.Item(15, 4).Value2 = 100
$ws.Cells.Item(15, 5).ClearContents()

$ws.Cells.Item(16, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Cells.Item(16, 2).Value2 = 34244
$ws.Cells.Item(16, 3).Value2 = ""
$ws.Cells.Item(16, 4).Value2 = 100
$ws.Cells.Item(16, 5).ClearContents()

$ws.Cells.Item(17, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Cells.Item(17, 2).Value2 = 442178
$ws.Cells.Item(17, 3).Value2 = ""
$ws.Cells.Item(17, 4).Value2 = 99.90000000000001

$ws.Cells.Item(18, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Cells.Item(18, 2).Value2 = 14239
$ws.Cells.Item(18, 3).Value2 = ""
$ws.Cells.Item(18, 4).Value2 = 100

$ws.Cells.Item(19, 1).Value2 = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Cells.Item(19, 2).Value2 = 265400
$ws.Cells.Item(19, 3).Value2 = ""
$ws.Cells.Item(19, 4).Value2 = 99.90000000000001

$ws.Cells.Item(20, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Cells.Item(20, 2).Value2 = 77849
$ws.Cells.Item(20, 3).Value2 = ""
$ws.Cells.Item(20, 4).Value2 = 99.90000000000001

$ws.Cells.Item(21, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Cells.Item(21, 2).Value2 = 59673
$ws.Cells.Item(21, 3).Value2 = ""
$ws.Cells.Item(21, 4).Value2 = 100

$ws.Cells.Item(22, 1).Value2 = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Cells.Item(22, 2).Value2 = 113652
$ws.Cells.Item(22, 3).Value2 = ""
$ws.Cells.Item(22, 4).Value2 = 100

# ---------------------------------------------------------------------
# 6) "Driver Vintage" column (E) holds literal yyyy-mm-dd text, not real
#    dates. Assigning a date-shaped string straight to .Value2 gets
#    auto-parsed into a date serial by the engine (just like typing one
#    into Excel), so force text mode via a "@" number format first, then
#    restore the plain right-aligned look (copied from a known-good data
#    cell) without disturbing the freshly-written string value.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

Set-TextValue $ws.Cells.Item(17, 5) "2024-11-10"
Set-TextValue $ws.Cells.Item(18, 5) "2022-05-23"
Set-TextValue $ws.Cells.Item(19, 5) "2022-05-01"
Set-TextValue $ws.Cells.Item(20, 5) "2021-08-18"
Set-TextValue $ws.Cells.Item(21, 5) "2020-08-05"
Set-TextValue $ws.Cells.Item(22, 5) "2019-12-14"

$ws.Cells.Item(3, 4).Copy()
$ws.Range("E17:E22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
